$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.179.81"
$ws.Range("E2").Value = "  -0.33%  "

$ws.Range("D3").Value = "3.052.56"
$ws.Range("E3").Value = "  +1.39%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'515.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.47%  "

$ws.Range("D6").Value = "'141.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.45%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.23%  "

$ws.Range("E8").Value = "  +1.38%  "

$ws.Range("D9").Value = "'7.23"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.40%  "

$ws.Range("D10").Value = "'0.109"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.46%  "

$ws.Range("E11").Value = "  +3.43%  "

$ws.Range("D12").Value = "3.570.82"
$ws.Range("E12").Value = "  +1.24%  "

$ws.Range("E13").Value = "  -3.00%  "

$ws.Range("D14").Value = "'27.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.68%  "

$ws.Range("D15").Value = "'0.0000166"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.22%  "

$ws.Range("D16").Value = "57.128.87"
$ws.Range("E16").Value = "  -0.49%  "

$ws.Range("E17").Value = "  -0.38%  "

$ws.Range("D18").Value = "3.049.94"
$ws.Range("E18").Value = "  +1.38%  "

$ws.Range("D19").Value = "'13.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.97%  "

$ws.Range("D20").Value = "'8.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.63%  "

$ws.Range("D21").Value = "'331.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.19%  "

$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("E23").Value = "  +1.95%  "

$ws.Range("D24").Value = "'65.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.89%  "

$ws.Range("D25").Value = "3.167.75"
$ws.Range("E25").Value = "  +1.00%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.36%  "

$ws.Range("D27").Value = "'0.165"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.93%  "

$ws.Range("D28").Value = "0.0₃0896"
$ws.Range("E28").Value = "  -2.19%  "

$ws.Range("D29").Value = "'6.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.36%  "

$ws.Range("D30").Value = "'7.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.59%  "

$ws.Range("E31").Value = "  +0.47%  "

$ws.Range("E32").Value = "  +1.62%  "

$ws.Range("D33").Value = "'20.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.34%  "

$ws.Range("E34").Value = "  -0.71%  "

$ws.Range("D35").Value = "'150.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.25%  "

$ws.Range("D36").Value = "'5.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.45%  "

$ws.Range("E37").Value = "  +0.45%  "

$ws.Range("D38").Value = "'25.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.00%  "

$ws.Range("E39").Value = "  +0.42%  "

$ws.Range("E40").Value = "  +1.20%  "

$ws.Range("E41").Value = "  -2.82%  "

$ws.Range("E42").Value = "  -0.12%  "

$ws.Range("D43").Value = "'0.663"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.17%  "

$ws.Range("E44").Value = "  -0.64%  "

$ws.Range("D45").Value = "2.205.99"
$ws.Range("E45").Value = "  -0.77%  "

$ws.Range("D46").Value = "'6.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.81%  "

$ws.Range("E47").Value = "  -2.74%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'20.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.65%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0242"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.61%  "

$ws.Range("E50").Value = "  +0.39%  "

$ws.Range("D51").Value = "'0.0172"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.19%  "
